# Refactor code to save results in a specified folder
# Update predicted IPC PO values (column C) for the sliding window results,
# then recompute the dependent DELTA (D), DELTA^2 (E) columns and the
# TOTAL / MSE summary rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newC = @(
    29.06111334257831,
    28.47071234469865,
    28.28848359845321,
    29.08335461496617,
    30.03117283040773,
    30.03955669991785,
    30.51717706393323,
    29.93871922786692,
    30.0226229808443,
    29.85382509620436,
    29.65558126129587,
    30.29504831114055,
    30.47089791408501,
    30.50014159004855,
    31.26975441792585,
    30.80844435579825,
    31.72773471084909,
    31.17942847998682,
    31.70892991584662,
    31.90208559243378,
    32.45973717409251,
    31.85781058353992,
    32.28776501006389,
    31.795646922987,
    32.78592873571744,
    32.75140148574432,
    32.28667916453107,
    33.77931662779437,
    32.6181679002533,
    32.92095187517251,
    33.56936814496778,
    33.8643056819737,
    34.52482003942425,
    34.44874248754216,
    35.08705449938859,
    35.12051217537368,
    35.47297676200688,
    35.67625688893011,
    36.07186189172979,
    36.74367428201054,
    38.30919573659418,
    38.66394966250266,
    38.88523822790722,
    39.16660625984799,
    39.91002835172897,
    39.79467645854881,
    40.14733291543243,
    40.46751621929954,
    41.47265485714376,
    40.6695322347838
)

$startRow = 2
$endRow = 51

$sumDelta = 0
$sumDeltaSq = 0

for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $startRow + $i
    $b = $ws.Cells.Item($row, 2).Value2
    $c = $newC[$i]
    $d = $c - $b
    $e = $d * $d

    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e

    $sumDelta = $sumDelta + $d
    $sumDeltaSq = $sumDeltaSq + $e
}

$count = $endRow - $startRow + 1
$mse = $sumDeltaSq / $count

# TOTAL row (52): C52 = sum of DELTA, E52 = sum of DELTA^2
$ws.Cells.Item(52, 3).Value = $sumDelta
$ws.Cells.Item(52, 5).Value = $sumDeltaSq

# MSE row (53): E53 = mean of DELTA^2
$ws.Cells.Item(53, 5).Value = $mse

Write-Output "Updated IPC PO predictions and recomputed deltas/summary rows."
